# Added sizing study for VA: switch r_w (row 7) distribution from
# "uniform" to "constant" on every K-scenario sheet, and leave the
# selection on D7 with the Virginia (iowa_k) sheet active.

$wb = $excel.ActiveWorkbook

$sheetNames = @("low_k", "med_low_k", "med_high_k", "high_k", "iowa_k")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("D7").Value = "constant"
    $ws.Range("D7").Select()
}

# Make the last sheet (iowa_k) the active tab, matching activeTab="4"
# and tabSelected="1" moving from the first sheet to this one.
$wsActive = $wb.Worksheets.Item("iowa_k")
$wsActive.Activate()
$wsActive.Range("D7").Select()
